$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 76: Warding Off Temptation | Enchanted Hardsilver Ink
$ws.Range("H76").Value = 6569
$ws.Range("I76").Value = 8970.294
$ws.Range("J76").Value = 3167.1667
$ws.Range("K76").Value = 8970.294
$ws.Range("L76").Value = 3167.1667
$ws.Range("M76").Value = -8655.294
$ws.Range("N76").Value = -3797.1667

# Row 79: The Garden of Arcane Delights (L) | Enchanted Hardsilver Ink
$ws.Range("H79").Value = 6569
$ws.Range("I79").Value = 8970.294
$ws.Range("J79").Value = 3167.1667
$ws.Range("K79").Value = 8970.294
$ws.Range("L79").Value = 3167.1667
$ws.Range("M79").Value = -7878.294
$ws.Range("N79").Value = -5351.1667

$ws = $wb.Worksheets.Item("ARM")
# Row 3: Skillet Labor | Bronze Skillet
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()

# Row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws.Range("H61").Value = 77078240
$ws.Range("I61").Value = 91001304
$ws.Range("K61").Value = 91001304
$ws.Range("M61").Value = -91001092

# Row 63: Rivets Run through It | Mythrite Rivets
$ws.Range("H63").Value = 4336
$ws.Range("I63").Value = 3800
$ws.Range("J63").Value = 5006
$ws.Range("K63").Value = 3800
$ws.Range("L63").Value = 5006
$ws.Range("M63").Value = -3114
$ws.Range("N63").Value = -6378

# Row 66: A Riveting Revival (L) | Mythrite Rivets
$ws.Range("H66").Value = 4336
$ws.Range("I66").Value = 3800
$ws.Range("J66").Value = 5006
$ws.Range("K66").Value = 19000
$ws.Range("L66").Value = 25030
$ws.Range("M66").Value = -15568
$ws.Range("N66").Value = -31894

# Row 74: As the Bolt Flies | Titanium Nugget
$ws.Range("H74").Value = 10501287
$ws.Range("I74").Value = 16734793
$ws.Range("J74").Value = 112108.89
$ws.Range("K74").Value = 16734793
$ws.Range("L74").Value = 112108.89
$ws.Range("M74").Value = -16733919
$ws.Range("N74").Value = -113856.89

# Row 77: Heavy Metal Banned (L) | Titanium Nugget
$ws.Range("H77").Value = 10501287
$ws.Range("I77").Value = 16734793
$ws.Range("J77").Value = 112108.89
$ws.Range("K77").Value = 83673965
$ws.Range("L77").Value = 560544.45
$ws.Range("M77").Value = -83669597
$ws.Range("N77").Value = -569280.45

# Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Range("H132").Value = 77001.82000000001
$ws.Range("I132").Value = 49573.906
$ws.Range("K132").Value = 148721.718
$ws.Range("M132").Value = -146191.718

# Row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws.Range("H136").Value = 77078240
$ws.Range("I136").Value = 91001304
$ws.Range("K136").Value = 273003912
$ws.Range("M136").Value = -273001362

$ws = $wb.Worksheets.Item("BSM")
# Row 5: Axe Me Anything | Bronze War Axe
$ws.Range("H5").Value = 1200
$ws.Range("J5").Value = 1200
$ws.Range("L5").Value = 1200
$ws.Range("N5").Value = -1426

# Row 105: Ingot to Wing It | Molybdenum Ingot
$ws.Range("H105").Value = 145137.14
$ws.Range("I105").Value = 178581.17
$ws.Range("K105").Value = 178581.17
$ws.Range("M105").Value = -176834.17

# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value = 1940.409
$ws.Range("I134").Value = 2209.2144
$ws.Range("K134").Value = 6627.6432
$ws.Range("M134").Value = -4092.6432

$ws = $wb.Worksheets.Item("CRP")
# Row 2: In with the New | Bone Harpoon
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()

# Row 58: You Do the Heavy Lifting | Mahogany Lumber
$ws.Range("H58").Value = 18183650
$ws.Range("I58").Value = 25001438
$ws.Range("J58").Value = 2885.8667
$ws.Range("K58").Value = 25001438
$ws.Range("L58").Value = 2885.8667
$ws.Range("M58").Value = -25001235
$ws.Range("N58").Value = -3291.8667

# Row 132: Hull Lotta Damage | Ginseng Lumber
$ws.Range("H132").Value = 51994.5
$ws.Range("J132").Value = 126876
$ws.Range("L132").Value = 380628
$ws.Range("N132").Value = -385688

# Row 134: Wood You Be Quiet | Ceiba Lumber
$ws.Range("H134").Value = 30868.297
$ws.Range("I134").Value = 1869.8148
$ws.Range("J134").Value = 109164.2
$ws.Range("K134").Value = 5609.4444
$ws.Range("L134").Value = 327492.6
$ws.Range("M134").Value = -3074.4444
$ws.Range("N134").Value = -332562.6

# Row 136: Turali Quality | Dark Mahogany Lumber
$ws.Range("H136").Value = 18183650
$ws.Range("I136").Value = 25001438
$ws.Range("J136").Value = 2885.8667
$ws.Range("K136").Value = 75004314
$ws.Range("L136").Value = 8657.6001
$ws.Range("M136").Value = -75001764
$ws.Range("N136").Value = -13757.6001

$ws = $wb.Worksheets.Item("CUL")
# Row 64: The Aroma of Faith | Baked Onion Soup
$ws.Range("H64").Value = 3453.625
$ws.Range("I64").Value = 1548
$ws.Range("J64").Value = 3725.8572
$ws.Range("K64").Value = 4644
$ws.Range("L64").Value = 11177.5716
$ws.Range("M64").Value = -4374
$ws.Range("N64").Value = -11717.5716

# Row 67: Soup's On (L) | Baked Onion Soup
$ws.Range("H67").Value = 3453.625
$ws.Range("I67").Value = 1548
$ws.Range("J67").Value = 3725.8572
$ws.Range("K67").Value = 4644
$ws.Range("L67").Value = 11177.5716
$ws.Range("M67").Value = -3708
$ws.Range("N67").Value = -13049.5716

# Row 97: The Frier Never Lies | Cottonseed Oil
$ws.Range("H97").Value = 3183
$ws.Range("I97").Value = 3183
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 9549
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -9053

# Row 117: A Good Omen | Peppered Popotoes
$ws.Range("H117").Value = 5556176.5
$ws.Range("I117").Value = 664.6667
$ws.Range("J117").Value = 11111688
$ws.Range("K117").Value = 1994.0001
$ws.Range("L117").Value = 33335064
$ws.Range("M117").Value = 1447.9999
$ws.Range("N117").Value = -33341948

# Row 132: More Mezcal | Cooking Mezcal
$ws.Range("H132").Value = 1426.8572
$ws.Range("I132").Value = 1108.4445
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 9976.0005
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -7446.0005
$ws.Range("N132").Value = -23060

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell | Hardsilver Ingot
$ws.Range("H80").Value = 8730.1
$ws.Range("I80").Value = 27550
$ws.Range("J80").Value = 4025.125
$ws.Range("K80").Value = 27550
$ws.Range("L80").Value = 4025.125
$ws.Range("M80").Value = -26552
$ws.Range("N80").Value = -6021.125

# Row 83: With a Noise That Reaches Heaven (L) | Hardsilver Ingot
$ws.Range("H83").Value = 8730.1
$ws.Range("I83").Value = 27550
$ws.Range("J83").Value = 4025.125
$ws.Range("K83").Value = 137750
$ws.Range("L83").Value = 20125.625
$ws.Range("M83").Value = -132758
$ws.Range("N83").Value = -30109.625

# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 107133.9
$ws.Range("I132").Value = 78057.766
$ws.Range("J132").Value = 170132.17
$ws.Range("K132").Value = 234173.298
$ws.Range("L132").Value = 510396.51
$ws.Range("M132").Value = -231643.298
$ws.Range("N132").Value = -515456.51

$ws = $wb.Worksheets.Item("LTW")
# Row 132: Tenets of Tanning | Silver Lobo Leather
$ws.Range("H132").Value = 37366.285
$ws.Range("I132").Value = 1142.3043
$ws.Range("J132").Value = 203996.6
$ws.Range("K132").Value = 3426.9129
$ws.Range("L132").Value = 611989.8
$ws.Range("M132").Value = -896.9129000000003
$ws.Range("N132").Value = -617049.8

$ws = $wb.Worksheets.Item("WVR")
# Row 28: Doublet Jeopardy | Cotton Doublet Vest of Gathering
$ws.Range("H28").Value = 4999
$ws.Range("J28").Value = 4999
$ws.Range("L28").Value = 4999
$ws.Range("N28").Value = -5695

# Row 126: A Polished Purchase | Snow Linen
$ws.Range("H126").Value = 1523.6923
$ws.Range("I126").Value = 1363.3158
$ws.Range("J126").Value = 1959
$ws.Range("K126").Value = 4089.9474
$ws.Range("L126").Value = 5877
$ws.Range("M126").Value = -1619.9474
$ws.Range("N126").Value = -10817

# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 72781.46000000001
$ws.Range("I132").Value = 48434.43
$ws.Range("K132").Value = 145303.29
$ws.Range("M132").Value = -142773.29

# Row 136: Weaving the Envelope | Sarcenet Cloth
$ws.Range("H136").Value = 40183.332
$ws.Range("I136").Value = 21809.375
$ws.Range("J136").Value = 334166.66
$ws.Range("K136").Value = 65428.125
$ws.Range("L136").Value = 1002499.98
$ws.Range("M136").Value = -62878.125
$ws.Range("N136").Value = -1007599.98
Write-Host "Applied scheduled runner updates to 8 sheets (31 rows)."
